{"js": "// The document text content shuffled across several paragraphs while every\n// paragraph's style/position and every run's formatting stayed the same.\n// We therefore rewrite the *text* of specific runs/paragraphs in place\n// rather than moving paragraphs around.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- Simple paragraphs: exactly one run, plain Replace keeps its rPr ----\nconst singleRunReplacements = [\n  // index, old text (sanity check only), new text\n  [5, \"O Trabalho de Gradua\u00e7\u00e3o (TG) tem por objetivo a integra\u00e7\u00e3o, o aprofundamento e aplica\u00e7\u00e3o dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realiza\u00e7\u00e3o de tarefas que fazem parte do perfil de atua\u00e7\u00e3o profissional do engenheiro f\u00edsico.\",\n       \"Elabora\u00e7\u00e3o, com a orienta\u00e7\u00e3o de um professor supervisor, de uma proposta de projeto em tema ligado \u00e0 \u00e1rea de ci\u00eancia e tecnologia.\"],\n  [6, \"The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer.\",\n       \"Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology.\"],\n  [10, \"Elabora\u00e7\u00e3o, com a orienta\u00e7\u00e3o de um professor supervisor, de uma proposta de projeto em tema ligado \u00e0 \u00e1rea de ci\u00eancia e tecnologia.\",\n        \"O aluno deve apresentar a proposta de trabalho \u00e0 uma banca formada pelo respons\u00e1vel pela disciplina e professores ou profissionais da \u00e1rea.\"],\n  [11, \"Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology.\",\n        \"The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer.\"],\n  [13, \"O aluno deve procurar um professor ou profissional com forma\u00e7\u00e3o na \u00e1rea de engenharia ou \u00e1reas correlatas, para a elabora\u00e7\u00e3o de uma proposta de projeto contendo motiva\u00e7\u00e3o e objetivos, fundamenta\u00e7\u00e3o te\u00f3rica e cronograma de execu\u00e7\u00e3o. O projeto propriamente dito ser\u00e1 desenvolvido e defendido na disciplina Trabalho de Gradua\u00e7\u00e3o II.\",\n        \"Avalia\u00e7\u00e3o e atribui\u00e7\u00e3o de nota do Trabalho de Gradua\u00e7\u00e3o por uma comiss\u00e3o de professores.\"],\n  [18, \"A ser definida no plano de trabalho.\",\n        \"1176388 - Luiz Tadeu Fernandes Eleno\"],\n];\n\nfor (const [idx, oldText, newText] of singleRunReplacements) {\n  paragraphs.items[idx].insertText(newText, \"Replace\");\n}\nawait context.sync();\n\n// ---- Paragraph 8 (\"Docente(s) Respons\u00e1vel(eis)\" bullet list): two runs,\n// first run ends with a line break. Rebuild via OOXML so the run/break\n// layout is preserved exactly. ----\nconst p8Ooxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr><w:r><w:t>O Trabalho de Gradua\u00e7\u00e3o (TG) tem por objetivo a integra\u00e7\u00e3o, o aprofundamento e aplica\u00e7\u00e3o dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realiza\u00e7\u00e3o de tarefas que fazem parte do perfil de atua\u00e7\u00e3o profissional do engenheiro f\u00edsico.</w:t><w:br/></w:r><w:r><w:t>O aluno deve procurar um professor ou profissional com forma\u00e7\u00e3o na \u00e1rea de engenharia ou \u00e1reas correlatas, para a elabora\u00e7\u00e3o de uma proposta de projeto contendo motiva\u00e7\u00e3o e objetivos, fundamenta\u00e7\u00e3o te\u00f3rica e cronograma de execu\u00e7\u00e3o. O projeto propriamente dito ser\u00e1 desenvolvido e defendido na disciplina Trabalho de Gradua\u00e7\u00e3o II.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nparagraphs.items[8].insertOoxml(p8Ooxml, \"Replace\");\nawait context.sync();\n\n// ---- Paragraph 16 (\"Avalia\u00e7\u00e3o\" bullet list): M\u00e9todo/Crit\u00e9rio/Norma de\n// recupera\u00e7\u00e3o bold labels stay put; only the three value runs change. ----\nconst p16Ooxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">M\u00e9todo: </w:t></w:r><w:r><w:t>A crit\u00e9rio da banca de avalia\u00e7\u00e3o poder\u00e1 ser estabelecido um prazo para readequa\u00e7\u00e3o e reapresenta\u00e7\u00e3o do plano de trabalho.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Crit\u00e9rio: </w:t></w:r><w:r><w:t>A ser definida no plano de trabalho.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Norma de recupera\u00e7\u00e3o: </w:t></w:r><w:r><w:t>5840730 - Antonio Jefferson da Silva Machado</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nparagraphs.items[16].insertOoxml(p16Ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# The document text content shuffled across several paragraphs while every\n# paragraph's style/position and every run's formatting stayed the same.\n# We therefore rewrite the *text* of specific runs/paragraphs in place\n# rather than moving paragraphs around. Word COM Paragraphs are 1-based.\n\n$d = $word.ActiveDocument\n\n# ---- Simple paragraphs: exactly one run, plain Range.Text keeps its rPr ----\n# Trim the trailing paragraph mark (last char of Range) before assigning so\n# the paragraph mark / pilcrow formatting is untouched.\nfunction Set-ParaText($doc, [int]$oneBasedIndex, [string]$newText) {\n    $p = $doc.Paragraphs.Item($oneBasedIndex)\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = $newText\n}\n\nSet-ParaText $d 6  \"Elabora\u00e7\u00e3o, com a orienta\u00e7\u00e3o de um professor supervisor, de uma proposta de projeto em tema ligado \u00e0 \u00e1rea de ci\u00eancia e tecnologia.\"\nSet-ParaText $d 7  \"Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology.\"\nSet-ParaText $d 11 \"O aluno deve apresentar a proposta de trabalho \u00e0 uma banca formada pelo respons\u00e1vel pela disciplina e professores ou profissionais da \u00e1rea.\"\nSet-ParaText $d 12 \"The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer.\"\nSet-ParaText $d 14 \"Avalia\u00e7\u00e3o e atribui\u00e7\u00e3o de nota do Trabalho de Gradua\u00e7\u00e3o por uma comiss\u00e3o de professores.\"\nSet-ParaText $d 19 \"1176388 - Luiz Tadeu Fernandes Eleno\"\n\n# ---- Paragraph 9 (\"Docente(s) Respons\u00e1vel(eis)\" bullet list): two runs,\n# first run ends with a line break. Rebuild via WordOpenXML so the run/break\n# layout is preserved exactly. ----\n$p9 = $d.Paragraphs.Item(9)\n$p9Ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr><w:r><w:t>O Trabalho de Gradua\u00e7\u00e3o (TG) tem por objetivo a integra\u00e7\u00e3o, o aprofundamento e aplica\u00e7\u00e3o dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realiza\u00e7\u00e3o de tarefas que fazem parte do perfil de atua\u00e7\u00e3o profissional do engenheiro f\u00edsico.</w:t><w:br/></w:r><w:r><w:t>O aluno deve procurar um professor ou profissional com forma\u00e7\u00e3o na \u00e1rea de engenharia ou \u00e1reas correlatas, para a elabora\u00e7\u00e3o de uma proposta de projeto contendo motiva\u00e7\u00e3o e objetivos, fundamenta\u00e7\u00e3o te\u00f3rica e cronograma de execu\u00e7\u00e3o. O projeto propriamente dito ser\u00e1 desenvolvido e defendido na disciplina Trabalho de Gradua\u00e7\u00e3o II.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$p9.Range.InsertXML($p9Ooxml)\n\n# ---- Paragraph 17 (\"Avalia\u00e7\u00e3o\" bullet list): M\u00e9todo/Crit\u00e9rio/Norma de\n# recupera\u00e7\u00e3o bold labels stay put; only the three value runs change. ----\n$p17 = $d.Paragraphs.Item(17)\n$p17Ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">M\u00e9todo: </w:t></w:r><w:r><w:t>A crit\u00e9rio da banca de avalia\u00e7\u00e3o poder\u00e1 ser estabelecido um prazo para readequa\u00e7\u00e3o e reapresenta\u00e7\u00e3o do plano de trabalho.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Crit\u00e9rio: </w:t></w:r><w:r><w:t>A ser definida no plano de trabalho.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Norma de recupera\u00e7\u00e3o: </w:t></w:r><w:r><w:t>5840730 - Antonio Jefferson da Silva Machado</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$p17.Range.InsertXML($p17Ooxml)\n"}
